$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (D) cells so numeric-looking strings
# are not auto-converted to numbers by Excel, matching the original
# inline-string (text) cell type.
$priceCells = @("D2", "D3", "D5", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D39", "D40", "D41", "D42", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($c in $priceCells) { $ws.Range($c).NumberFormat = "@" }

# Apply the updated cell values
$ws.Range("D2").Value = "29.426.24"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").Value = "1.899.82"
$ws.Range("E3").Value = "  -0.57%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "325.72"
$ws.Range("E5").Value = "  -1.89%  "
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").Value = "0.4829"
$ws.Range("E7").Value = "  +3.38%  "
$ws.Range("D8").Value = "0.4065"
$ws.Range("E8").Value = "  -0.92%  "
$ws.Range("D9").Value = "0.08075"
$ws.Range("E9").Value = "  +0.61%  "
$ws.Range("D10").Value = "1.003"
$ws.Range("E10").Value = "  -1.02%  "
$ws.Range("D11").Value = "23.50"
$ws.Range("E11").Value = "  +5.15%  "
$ws.Range("D12").Value = "1.898.11"
$ws.Range("E12").Value = "  -2.08%  "
$ws.Range("D13").Value = "5.970"
$ws.Range("E13").Value = "  -0.08%  "
$ws.Range("D14").Value = "7.074"
$ws.Range("E14").Value = "  -1.48%  "
$ws.Range("D15").Value = "90.07"
$ws.Range("E15").Value = "  +0.40%  "
$ws.Range("E16").Value = "  +0.31%  "
$ws.Range("E17").Value = "  +1.63%  "
$ws.Range("D18").Value = "0.00001032"
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("D19").Value = "17.63"
$ws.Range("E19").Value = "  -0.58%  "
$ws.Range("D20").Value = "1.003"
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("D21").Value = "29.448.85"
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").Value = "5.547"
$ws.Range("E22").Value = "  -0.36%  "
$ws.Range("D23").Value = "11.77"
$ws.Range("E23").Value = "  +1.94%  "
$ws.Range("D24").Value = "2.154"
$ws.Range("E24").Value = "  -2.74%  "
$ws.Range("D25").Value = "2.134.81"
$ws.Range("E25").Value = "  -0.79%  "
$ws.Range("D26").Value = "153.78"
$ws.Range("E26").Value = "  -0.21%  "
$ws.Range("D27").Value = "19.90"
$ws.Range("E27").Value = "  +0.34%  "
$ws.Range("D28").Value = "6.084"
$ws.Range("E28").Value = "  +5.79%  "
$ws.Range("E29").Value = "  -2.30%  "
$ws.Range("D30").Value = "118.58"
$ws.Range("E30").Value = "  +1.18%  "
$ws.Range("D31").Value = "1.035"
$ws.Range("E31").Value = "  -2.56%  "
$ws.Range("D32").Value = "0.09512"
$ws.Range("E32").Value = "  +0.60%  "
$ws.Range("D33").Value = "5.518"
$ws.Range("E33").Value = "  +2.21%  "
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "1.390"
$ws.Range("E34").Value = "  -2.44%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "3.543"
$ws.Range("E35").Value = "  -1.01%  "
$ws.Range("D36").Value = "0.06083"
$ws.Range("E36").Value = "  -0.70%  "
$ws.Range("E37").Value = "  -0.35%  "
$ws.Range("E38").Value = "  -0.69%  "
$ws.Range("D39").Value = "0.5893"
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("D40").Value = "7.906"
$ws.Range("E40").Value = "  -5.96%  "
$ws.Range("D41").Value = "0.1846"
$ws.Range("E41").Value = "  +0.28%  "
$ws.Range("D42").Value = "10.22"
$ws.Range("E42").Value = "  +0.34%  "
$ws.Range("E43").Value = "  +0.62%  "
$ws.Range("D44").Value = "0.07852"
$ws.Range("E44").Value = "  +4.54%  "
$ws.Range("D45").Value = "2.384"
$ws.Range("E45").Value = "  +2.11%  "
$ws.Range("D46").Value = "12.30"
$ws.Range("E46").Value = "  +1.12%  "
$ws.Range("D47").Value = "0.5537"
$ws.Range("E47").Value = "  -0.38%  "
$ws.Range("D48").Value = "1.925"
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("D49").Value = "114.30"
$ws.Range("E49").Value = "  +1.10%  "
$ws.Range("D50").Value = "72.36"
$ws.Range("E50").Value = "  +1.19%  "
$ws.Range("D51").Value = "0.2928"
$ws.Range("E51").Value = "  -1.39%  "

# Restore the default "Normal" style on the Price cells so only the
# value changes (not the cell style/format) remain in the diff.
foreach ($c in $priceCells) { $ws.Range($c).Style = "Normal" }
